$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "''Akkermansia_muciniphila_ATCC_BAA_835.mat'"
$ws.Range("C2").Value = 0

$ws.Range("B3").Value = "''Alistipes_finegoldii_DSM_17242.mat'"
$ws.Range("C3").Value = 0

$ws.Range("B4").Value = "''Alistipes_putredinis_DSM_17216.mat'"
$ws.Range("C4").Value = 0

$ws.Range("B5").Value = "''Bacteroides_cellulosilyticus_DSM_14838.mat'"
$ws.Range("C5").Value = 0

$ws.Range("B6").Value = "''Bacteroides_coprophilus_DSM_18228.mat'"
$ws.Range("C6").Value = 0

$ws.Range("B7").Value = "''Bacteroides_fragilis_3_1_12.mat'"
$ws.Range("C7").Value = 0

$ws.Range("B8").Value = "''Bacteroides_oleiciplenus_YIT_12058.mat'"
$ws.Range("C8").Value = 0

$ws.Range("B9").Value = "''Bacteroides_ovatus_ATCC_8483.mat'"
$ws.Range("C9").Value = 0

$ws.Range("B10").Value = "''Bacteroides_plebeius_M12_DSM_17135.mat'"
$ws.Range("C10").Value = 0

$ws.Range("B11").Value = "''Bacteroides_salyersiae_WAL_10018.mat'"
$ws.Range("C11").Value = 0

$ws.Range("B12").Value = "''Bacteroides_stercoris_ATCC_43183.mat'"
$ws.Range("C12").Value = 0

$ws.Range("B13").Value = "''Bacteroides_thetaiotaomicron_VPI_5482.mat'"
$ws.Range("C13").Value = 0

$ws.Range("B14").Value = "''Bacteroides_uniformis_ATCC_8492.mat'"
$ws.Range("C14").Value = 0

$ws.Range("B15").Value = "''Bacteroides_vulgatus_ATCC_8482.mat'"
$ws.Range("C15").Value = 0

$ws.Range("B16").Value = "''Barnesiella_intestinihominis_YIT_11860.mat'"
$ws.Range("C16").Value = 0

$ws.Range("B17").Value = "''Bifidobacterium_animalis_lactis_AD011.mat'"
$ws.Range("C17").Value = 0

$ws.Range("B18").Value = "''Bilophila_wadsworthia_3_1_6.mat'"
$ws.Range("C18").Value = 0

$ws.Range("B19").Value = "''Dorea_longicatena_DSM_13814.mat'"
$ws.Range("C19").Value = 0

$ws.Range("B20").Value = "''Escherichia_coli_O157_H7_str_Sakai.mat'"
$ws.Range("C20").Value = 0.006

$ws.Range("B21").Value = "''Eubacterium_limosum_KIST612.mat'"
$ws.Range("C21").Value = 0

$ws.Range("B22").Value = "''Eubacterium_ramulus_ATCC_29099.mat'"
$ws.Range("C22").Value = 0

$ws.Range("B23").Value = "''Eubacterium_ventriosum_ATCC_27560.mat'"
$ws.Range("C23").Value = 0

$ws.Range("B24").Value = "''Flavonifractor_plautii_ATCC_29863.mat'"
$ws.Range("C24").Value = 0

$ws.Range("B25").Value = "''Marvinbryantia_formatexigens_I_52_DSM_14469.mat'"
$ws.Range("C25").Value = 0

$ws.Range("B26").Value = "''Odoribacter_splanchnicus_1651_6_DSM_20712.mat'"
$ws.Range("C26").Value = 0

$ws.Range("B27").Value = "''Parabacteroides_distasonis_ATCC_8503.mat'"
$ws.Range("C27").Value = 0

$ws.Range("B28").Value = "''Parabacteroides_johnsonii_DSM_18315.mat'"
$ws.Range("C28").Value = 0.042

$ws.Range("B29").Value = "''Paraprevotella_xylaniphila_YIT_11841.mat'"
$ws.Range("C29").Value = 0

$ws.Range("B30").Value = "''Parasutterella_excrementihominis_YIT_11859.mat'"
$ws.Range("C30").Value = 0

$ws.Range("B31").Value = "''Phascolarctobacterium_succinatutens_YIT_12067.mat'"
$ws.Range("C31").Value = 0.046

$ws.Range("B32").Value = "''Prevotella_copri_CB7_DSM_18205.mat'"
$ws.Range("C32").Value = 0

$ws.Range("B33").Value = "''Prevotella_stercorea_DSM_18206.mat'"
$ws.Range("C33").Value = 0.906

$ws.Range("B34").Value = "''Roseburia_inulinivorans_DSM_16841.mat'"
$ws.Range("C34").Value = 0

$ws.Range("B35").Value = "''Sutterella_wadsworthensis_3_1_45B.mat'"
$ws.Range("C35").Value = 0
